# Updates cryptos list values (price/volume) and row reorder for TRON/Polkadot and Toncoin/PancakeSwap
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.672.94"
$ws.Range("E2").Value = "  +2.90%  "
$ws.Range("D3").Value = "2.198.70"
$ws.Range("E3").Value = "  +1.69%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'250.82"
$ws.Range("E5").Value = "  +5.75%  "
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("D7").Value = "'74.66"
$ws.Range("E7").Value = "  +5.10%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.588"
$ws.Range("E9").Value = "  +2.67%  "
$ws.Range("D10").Value = "'40.28"
$ws.Range("E10").Value = "  +1.77%  "
$ws.Range("D11").Value = "'0.0920"
$ws.Range("E11").Value = "  +2.43%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.101"
$ws.Range("E12").Value = "  +1.59%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'6.80"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").Value = "2.530.23"
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("D15").Value = "'14.34"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "2.209.78"
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("D17").Value = "'0.781"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "42.565.70"
$ws.Range("E18").Value = "  +3.00%  "
$ws.Range("D19").Value = "'0.0000102"
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("D20").Value = "'71.17"
$ws.Range("E20").Value = "  +2.44%  "
$ws.Range("E21").Value = "  +3.08%  "
$ws.Range("D22").Value = "'228.80"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("E23").Value = "  +9.17%  "
$ws.Range("D24").Value = "'9.48"
$ws.Range("E24").Value = "  -4.83%  "
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").Value = "'10.71"
$ws.Range("E26").Value = "  +0.53%  "
$ws.Range("E27").Value = "  +3.73%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'2.20"
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.22"
$ws.Range("E29").Value = "  +2.54%  "
$ws.Range("D30").Value = "'37.61"
$ws.Range("E30").Value = "  +13.85%  "
$ws.Range("D31").Value = "'169.42"
$ws.Range("E31").Value = "  -1.09%  "
$ws.Range("D32").Value = "'20.10"
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("D33").Value = "'0.0799"
$ws.Range("E33").Value = "  +3.66%  "
$ws.Range("E34").Value = "  +1.67%  "
$ws.Range("E35").Value = "  +0.96%  "
$ws.Range("E36").Value = "  +4.06%  "
$ws.Range("E37").Value = "  +4.01%  "
$ws.Range("D38").Value = "'0.0327"
$ws.Range("E38").Value = "  +9.18%  "
$ws.Range("D39").Value = "'12.25"
$ws.Range("E39").Value = "  +2.40%  "
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("E41").Value = "  +5.79%  "
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("D43").Value = "'59.18"
$ws.Range("E43").Value = "  +1.09%  "
$ws.Range("D44").Value = "'103.09"
$ws.Range("E44").Value = "  +7.41%  "
$ws.Range("D45").Value = "'0.480"
$ws.Range("E45").Value = "  +23.84%  "
$ws.Range("D46").Value = "'8.47"
$ws.Range("E46").Value = "  +1.49%  "
$ws.Range("E47").Value = "  +2.51%  "
$ws.Range("D48").Value = "'2.42"
$ws.Range("E48").Value = "  +12.75%  "
$ws.Range("E49").Value = "  +2.64%  "
$ws.Range("E50").Value = "  +2.03%  "

Write-Host "Applied cryptos update."
